$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 1350
$ws.Range("I16").Value = 1475
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 1475
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -1245
$ws.Range("N16").Value = -1560

$ws.Range("H132").Value = 4535.8906
$ws.Range("I132").Value = 3460.623
$ws.Range("J132").Value = 10001.833
$ws.Range("K132").Value = 10381.869
$ws.Range("L132").Value = 30005.499
$ws.Range("M132").Value = -7851.869000000001
$ws.Range("N132").Value = -35065.499

$ws.Range("H135").Value = 980.5
$ws.Range("I135").Value = 667.9722
$ws.Range("J135").Value = 2386.875
$ws.Range("K135").Value = 6011.749800000001
$ws.Range("L135").Value = 21481.875
$ws.Range("M135").Value = -3476.749800000001
$ws.Range("N135").Value = -26551.875

$ws.Range("H137").Value = 2049.9534
$ws.Range("I137").Value = 2686.9412
$ws.Range("K137").Value = 8060.823600000001
$ws.Range("M137").Value = -5510.823600000001

$ws.Range("H138").Value = 2702.4795
$ws.Range("I138").Value = 1657.0667
$ws.Range("J138").Value = 3431.8372
$ws.Range("K138").Value = 4971.2001
$ws.Range("L138").Value = 10295.5116
$ws.Range("M138").Value = 168.7999
$ws.Range("N138").Value = -20575.5116

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5095.69
$ws.Range("I32").Value = 5095.69
$ws.Range("K32").Value = 5095.69
$ws.Range("M32").Value = -4808.69

$ws.Range("H102").Value = 2460.625
$ws.Range("I102").Value = 1754.3334
$ws.Range("J102").Value = 3368.7144
$ws.Range("K102").Value = 1754.3334
$ws.Range("L102").Value = 3368.7144
$ws.Range("M102").Value = -132.3334
$ws.Range("N102").Value = -6612.7144

$ws.Range("H122").Value = 1604.3889
$ws.Range("I122").Value = 1871.4286
$ws.Range("J122").Value = 1434.4546
$ws.Range("K122").Value = 5614.2858
$ws.Range("L122").Value = 4303.3638
$ws.Range("M122").Value = -3164.2858
$ws.Range("N122").Value = -9203.363799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1455.8823
$ws.Range("I22").Value = 930
$ws.Range("J22").Value = 5400
$ws.Range("K22").Value = 930
$ws.Range("L22").Value = 5400
$ws.Range("M22").Value = -757
$ws.Range("N22").Value = -5746

$ws.Range("H75").Value = 7927.5454
$ws.Range("I75").Value = 4535.6665
$ws.Range("K75").Value = 4535.6665
$ws.Range("M75").Value = -3599.6665

$ws.Range("H78").Value = 7927.5454
$ws.Range("I78").Value = 4535.6665
$ws.Range("K78").Value = 13606.9995
$ws.Range("M78").Value = -8926.999500000002

$ws.Range("H86").Value = 1184.5714
$ws.Range("I86").Value = 1184.5714
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1184.5714
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -61.57140000000004
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 1184.5714
$ws.Range("I89").Value = 1184.5714
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5922.857
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -306.857
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1831.3062
$ws.Range("I132").Value = 1291.6428
$ws.Range("J132").Value = 2550.8572
$ws.Range("K132").Value = 3874.9284
$ws.Range("L132").Value = 7652.571599999999
$ws.Range("M132").Value = -1344.9284
$ws.Range("N132").Value = -12712.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5924.1177
$ws.Range("I3").Value = 1673
$ws.Range("J3").Value = 11997.143
$ws.Range("K3").Value = 5019
$ws.Range("L3").Value = 35991.429
$ws.Range("M3").Value = -4907
$ws.Range("N3").Value = -36215.429

$ws.Range("H5").Value = 664.2353000000001
$ws.Range("I5").Value = 394.46155
$ws.Range("J5").Value = 1541
$ws.Range("K5").Value = 1183.38465
$ws.Range("L5").Value = 4623
$ws.Range("M5").Value = -1071.38465
$ws.Range("N5").Value = -4847

$ws.Range("H111").Value = 2010
$ws.Range("I111").Value = 500
$ws.Range("K111").Value = 1500
$ws.Range("M111").Value = 1567

$ws.Range("H126").Value = 2562.5
$ws.Range("I126").Value = 1838.3334
$ws.Range("J126").Value = 2997
$ws.Range("K126").Value = 5515.0002
$ws.Range("L126").Value = 8991
$ws.Range("M126").Value = -575.0002000000004
$ws.Range("N126").Value = -18871

$ws.Range("H135").Value = 664.2353000000001
$ws.Range("I135").Value = 394.46155
$ws.Range("J135").Value = 1541
$ws.Range("K135").Value = 3550.15395
$ws.Range("L135").Value = 13869
$ws.Range("M135").Value = -1015.15395
$ws.Range("N135").Value = -18939

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 32.6
$ws.Range("I2").Value = 57.142857
$ws.Range("J2").Value = 25.130434
$ws.Range("K2").Value = 57.142857
$ws.Range("L2").Value = 25.130434
$ws.Range("M2").Value = 55.857143
$ws.Range("N2").Value = -251.130434

$ws.Range("H57").Value = 2949
$ws.Range("I57").Value = 2949
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 2949
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -2129
$ws.Range("N57").ClearContents()

$ws.Range("H80").Value = 3275.75
$ws.Range("I80").Value = 2925
$ws.Range("J80").Value = 3626.5
$ws.Range("K80").Value = 2925
$ws.Range("L80").Value = 3626.5
$ws.Range("M80").Value = -1927
$ws.Range("N80").Value = -5622.5

$ws.Range("H83").Value = 3275.75
$ws.Range("I83").Value = 2925
$ws.Range("J83").Value = 3626.5
$ws.Range("K83").Value = 14625
$ws.Range("L83").Value = 18132.5
$ws.Range("M83").Value = -9633
$ws.Range("N83").Value = -28116.5

$ws.Range("H122").Value = 3234.7334
$ws.Range("I122").Value = 2972.1
$ws.Range("J122").Value = 3760
$ws.Range("K122").Value = 8916.299999999999
$ws.Range("L122").Value = 11280
$ws.Range("M122").Value = -6466.299999999999
$ws.Range("N122").Value = -16180

$ws.Range("H126").Value = 5210861.5
$ws.Range("I126").Value = 9617594
$ws.Range("J126").Value = 2904.9092
$ws.Range("K126").Value = 28852782
$ws.Range("L126").Value = 8714.7276
$ws.Range("M126").Value = -28850312
$ws.Range("N126").Value = -13654.7276

$ws.Range("H132").Value = 3321.862
$ws.Range("I132").Value = 2225.2307
$ws.Range("J132").Value = 4212.875
$ws.Range("K132").Value = 6675.6921
$ws.Range("L132").Value = 12638.625
$ws.Range("M132").Value = -4145.6921
$ws.Range("N132").Value = -17698.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 4300
$ws.Range("J41").Value = 4300
$ws.Range("L41").Value = 4300
$ws.Range("N41").Value = -5080

$ws.Range("H45").Value = 8558
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 8558
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 8558
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -9540
